$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 719; this shifts the existing rows 719:748
# down to 720:749 and extends the sheet dimension to R749.
$ws.Rows("719:719").Insert()

# Populate the newly inserted row 719 with the new weekly record.
$ws.Range("A719").Value = 3
$ws.Range("B719").Value = "Femacal de La Calera"
$ws.Range("C719").Value = "Coquimbo"
$ws.Range("D719").Value = 45147
$ws.Range("E719").Value = 5
$ws.Range("F719").Value = 100112021
$ws.Range("G719").Value = "Ají"
$ws.Range("H719").Value = "Inferno"
$ws.Range("I719").Value = "Primera"
$ws.Range("J719").Value = 73
$ws.Range("K719").Value = 15000
$ws.Range("L719").Value = 16000
$ws.Range("M719").Value = 15521
$ws.Range("N719").Value = "$/caja 10 kilos"
$ws.Range("O719").Value = "Región de Arica y Parinacota"
$ws.Range("P719").Value = 1552
$ws.Range("Q719").Value = 10
$ws.Range("R719").Value = "Hortaliza"
